$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the new row 14 values
$ws.Range("A14").Value = "Exp 18"
$ws.Range("B14").Value = 0.5
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "Local"
$ws.Range("E14").Value = -1
$ws.Range("F14").Value = "Exp 18.png"

# Copy formatting from the row above (A13:E13) onto the new row's cells (A14:E14)
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)

# Move the active selection to B14 to match the saved workbook state
$ws.Range("B14").Select()
